$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 213, pushing existing rows 213-216 down to 215-218.
$ws.Range("A213:T214").Insert()

# Populate the first new row (213)
$ws.Cells.Item(213, 1).Value = 4
$ws.Cells.Item(213, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(213, 3).Value = "Los Lagos"
$ws.Cells.Item(213, 4).Value = 44595
$ws.Cells.Item(213, 5).Value = 10
$ws.Cells.Item(213, 6).Value = "Fruta"
$ws.Cells.Item(213, 7).Value = 100102
$ws.Cells.Item(213, 8).Value = "Cítricos"
$ws.Cells.Item(213, 9).Value = 100102006
$ws.Cells.Item(213, 10).Value = "Pomelo"
$ws.Cells.Item(213, 11).Value = "Start Ruby"
$ws.Cells.Item(213, 12).Value = "Primera"
$ws.Cells.Item(213, 13).Value = 80
$ws.Cells.Item(213, 14).Value = 13000
$ws.Cells.Item(213, 15).Value = 14000
$ws.Cells.Item(213, 16).Value = 13500
$ws.Cells.Item(213, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(213, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(213, 19).Value = 964
$ws.Cells.Item(213, 20).Value = 14

# Populate the second new row (214)
$ws.Cells.Item(214, 1).Value = 4
$ws.Cells.Item(214, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(214, 3).Value = "Los Lagos"
$ws.Cells.Item(214, 4).Value = 44595
$ws.Cells.Item(214, 5).Value = 10
$ws.Cells.Item(214, 6).Value = "Fruta"
$ws.Cells.Item(214, 7).Value = 100102
$ws.Cells.Item(214, 8).Value = "Cítricos"
$ws.Cells.Item(214, 9).Value = 100102006
$ws.Cells.Item(214, 10).Value = "Pomelo"
$ws.Cells.Item(214, 11).Value = "Start Ruby"
$ws.Cells.Item(214, 12).Value = "Segunda"
$ws.Cells.Item(214, 13).Value = 40
$ws.Cells.Item(214, 14).Value = 11000
$ws.Cells.Item(214, 15).Value = 11000
$ws.Cells.Item(214, 16).Value = 11000
$ws.Cells.Item(214, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(214, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(214, 19).Value = 786
$ws.Cells.Item(214, 20).Value = 14
